$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = 42619.893796296295
$ws.Cells.Item($row, 2).Value = 28
$ws.Cells.Item($row, 3).Value = 63
$ws.Cells.Item($row, 4).Value = 33
$ws.Cells.Item($row, 5).Value = 63
$ws.Cells.Item($row, 6).Value = 33
$ws.Cells.Item($row, 7).Value = 27812
$ws.Cells.Item($row, 8).Value = 14564
$ws.Cells.Item($row, 9).Value = 2284
$ws.Cells.Item($row, 10).Value = 355
$ws.Cells.Item($row, 11).Value = 184
$ws.Cells.Item($row, 12).Value = 18
$ws.Cells.Item($row, 13).Value = 9
$ws.Cells.Item($row, 14).Value = "Noun"
